# Populate the new "carrier" (D) and "pair_kind" (J) columns, plus the new
# unique_video / unique_audio rows (C/D) in rows 14-21, that were added to
# the stimuli sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("carrier") for the practice block (rows 2-5) ---
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# --- Column J ("pair_kind") for the generic block (rows 6-9) ---
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# --- New rows 14-21: unique_video / unique_audio kind + carrier word ---
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "can"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "can"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "do"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "do"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "look"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "look"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "where"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "where"

Write-Output "applied stimuli updates"
